$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "case with 380 kV done"
# The slack-bus voltage setpoint was changed from 1.05 pu to 1.02 pu, and the
# corresponding bus voltage-magnitude results (res_bus/vm_pu) were recomputed for
# every row (buses 0-23 in column A, stored in sheet rows 2-25). Column G (slack bus,
# always 1 pu) and the empty column H are unaffected by the change.

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.015334624463914
$ws.Range("D2").Value = 1.022227495211961
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.027398919303818
$ws.Range("I2").Value = 1.027797478233221
$ws.Range("J2").Value = 1.020560732494065
$ws.Range("K2").Value = 1.025062496934452
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.030218762313519
$ws.Range("N2").Value = 1.01097423641684

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.016157263640075
$ws.Range("D3").Value = 1.022823782703563
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.028465428481951
$ws.Range("I3").Value = 1.027912724250136
$ws.Range("J3").Value = 1.02101901036786
$ws.Range("K3").Value = 1.025465980798706
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.03109228877894
$ws.Range("N3").Value = 1.011126600144823

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.016689968287298
$ws.Range("D4").Value = 1.023209757317384
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.029156137189782
$ws.Range("I4").Value = 1.027986004275579
$ws.Range("J4").Value = 1.021315317404754
$ws.Range("K4").Value = 1.025726503443491
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.031657533419213
$ws.Range("N4").Value = 1.01122508176115

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.016914012252309
$ws.Range("D5").Value = 1.02337205217865
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.029446655390641
$ws.Range("I5").Value = 1.028016501358796
$ws.Range("J5").Value = 1.02143982896079
$ws.Range("K5").Value = 1.025835892373616
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.031895164929443
$ws.Range("N5").Value = 1.011266457231374

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.016951635748051
$ws.Range("D6").Value = 1.023399303979053
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.029495443111034
$ws.Range("I6").Value = 1.028021603767645
$ws.Range("J6").Value = 1.021460731690366
$ws.Range("K6").Value = 1.025854251319912
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.031935064437671
$ws.Range("N6").Value = 1.011273402807694

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.016692961602371
$ws.Range("D7").Value = 1.023211925788948
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.029160018542362
$ws.Range("I7").Value = 1.027986412997262
$ws.Range("J7").Value = 1.021316981355259
$ws.Range("K7").Value = 1.025727965634756
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.031660708653659
$ws.Range("N7").Value = 1.011225634725625

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.015612555169182
$ws.Range("D8").Value = 1.022428984088647
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.027759225121192
$ws.Range("I8").Value = 1.02783669331013
$ws.Range("J8").Value = 1.020715656540081
$ws.Range("K8").Value = 1.025198971149002
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.030513971146352
$ws.Range("N8").Value = 1.011025750572758

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.013711879274531
$ws.Range("D9").Value = 1.021050458576415
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.025295524510665
$ws.Range("I9").Value = 1.027562999064489
$ws.Range("J9").Value = 1.019654332854543
$ws.Range("K9").Value = 1.02426258751559
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.028493421686888
$ws.Range("N9").Value = 1.010672718997262

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.012446950444736
$ws.Range("D10").Value = 1.020132288005091
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.023656247546513
$ws.Range("I10").Value = 1.027373934387244
$ws.Range("J10").Value = 1.018945691653018
$ws.Range("K10").Value = 1.023635554809045
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.027146536489323
$ws.Range("N10").Value = 1.010436841680071

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.011899757658281
$ws.Range("D11").Value = 1.01973493071277
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.022947188261094
$ws.Range("I11").Value = 1.027290508410251
$ws.Range("J11").Value = 1.018638595532318
$ws.Range("K11").Value = 1.023363397779722
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.026563364970906
$ws.Range("N11").Value = 1.010334584709735

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.011696586764925
$ws.Range("D12").Value = 1.019587368667871
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.022683926679882
$ws.Range("I12").Value = 1.027259286494431
$ws.Range("J12").Value = 1.01852448980548
$ws.Range("K12").Value = 1.023262210312098
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.026346755617733
$ws.Range("N12").Value = 1.010296584168348

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.011740163931017
$ws.Range("D13").Value = 1.019619019656595
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.022740391985831
$ws.Range("I13").Value = 1.027265994272372
$ws.Range("J13").Value = 1.018548967509721
$ws.Range("K13").Value = 1.023283919704567
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.026393218753118
$ws.Range("N13").Value = 1.010304736210915

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01188296183168
$ws.Range("D14").Value = 1.019722732482303
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.022925424634868
$ws.Range("I14").Value = 1.027287932363154
$ws.Range("J14").Value = 1.018629164257395
$ws.Range("K14").Value = 1.023355035548861
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.026545459837211
$ws.Range("N14").Value = 1.01033144393314

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.011970955085847
$ws.Range("D15").Value = 1.01978663794962
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.023039444563272
$ws.Range("I15").Value = 1.02730141818112
$ws.Range("J15").Value = 1.018678571305861
$ws.Range("K15").Value = 1.023398839648771
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.026639261478127
$ws.Range("N15").Value = 1.010347897101013

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.012483277098438
$ws.Range("D16").Value = 1.020158664000661
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.023703321512
$ws.Range("I16").Value = 1.027379438278857
$ws.Range("J16").Value = 1.01896606743141
$ws.Range("K16").Value = 1.023653603420488
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.027185240537101
$ws.Range("N16").Value = 1.010443625629575

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.012804785894576
$ws.Range("D17").Value = 1.020392085241425
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.024119957506873
$ws.Range("I17").Value = 1.027427961021817
$ws.Range("J17").Value = 1.019146339965202
$ws.Range("K17").Value = 1.023813237354652
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.027527729624272
$ws.Range("N17").Value = 1.010503641589674

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.012992367562369
$ws.Range("D18").Value = 1.020528256695994
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.024363047369521
$ws.Range("I18").Value = 1.027456112990025
$ws.Range("J18").Value = 1.019251465716326
$ws.Range("K18").Value = 1.023906286558499
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.027727501396962
$ws.Range("N18").Value = 1.010538636229823

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.013056336679499
$ws.Range("D19").Value = 1.020574691121279
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.024445947115801
$ws.Range("I19").Value = 1.027465686535363
$ws.Range("J19").Value = 1.019287306751763
$ws.Range("K19").Value = 1.023938003286553
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.027795619012084
$ws.Range("N19").Value = 1.010550566510106

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.012770285756281
$ws.Range("D20").Value = 1.020367039188977
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.024075248838023
$ws.Range("I20").Value = 1.02742277055652
$ws.Range("J20").Value = 1.01912700092114
$ws.Range("K20").Value = 1.023796116608251
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.027490983400432
$ws.Range("N20").Value = 1.010497203647966

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.011840909166322
$ws.Range("D21").Value = 1.019692190684685
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.022870933962483
$ws.Range("I21").Value = 1.027281478591736
$ws.Range("J21").Value = 1.018605549315296
$ws.Range("K21").Value = 1.023334096352925
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.026500628446934
$ws.Range("N21").Value = 1.01032357966171

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.01125704175007
$ws.Range("D22").Value = 1.019268084991334
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.022114396818435
$ws.Range("I22").Value = 1.027191290335294
$ws.Range("J22").Value = 1.018277480985683
$ws.Range("K22").Value = 1.023043049664325
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.025877991023498
$ws.Range("N22").Value = 1.010214312854077

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.011566516074194
$ws.Range("D23").Value = 1.019492892109626
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.022515388367321
$ws.Range("I23").Value = 1.027239228847866
$ws.Range("J23").Value = 1.018451415890014
$ws.Range("K23").Value = 1.023197391381909
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.026208059059641
$ws.Range("N23").Value = 1.010272246853723

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.012785874726199
$ws.Range("D24").Value = 1.020378356357923
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.024095450535016
$ws.Range("I24").Value = 1.027425116369524
$ws.Range("J24").Value = 1.019135739477488
$ws.Range("K24").Value = 1.023803852930041
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.027507587425883
$ws.Range("N24").Value = 1.010500112712754

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.014202869572476
$ws.Range("D25").Value = 1.021406698204939
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.02593189127462
$ws.Range("I25").Value = 1.027634921663868
$ws.Range("J25").Value = 1.019928907292634
$ws.Range("K25").Value = 1.024505159378496
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.029015759625427
$ws.Range("N25").Value = 1.010764079883495
